$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 401; this shifts existing rows 401-421 down to 402-422
$ws.Rows.Item(401).Insert()

# Populate the newly inserted row 401 with the new record
$ws.Cells.Item(401, 1).Value = 9
$ws.Cells.Item(401, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(401, 3).Value = "Metropolitana"
$ws.Cells.Item(401, 4).Value = 45041
$ws.Cells.Item(401, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(401, 5).Value = 13
$ws.Cells.Item(401, 6).Value = 300000001
$ws.Cells.Item(401, 7).Value = "Rabanito"
$ws.Cells.Item(401, 8).Value = "Sin especificar"
$ws.Cells.Item(401, 9).Value = "Primera"
$ws.Cells.Item(401, 10).Value = 7000
$ws.Cells.Item(401, 11).Value = 3000
$ws.Cells.Item(401, 12).Value = 3000
$ws.Cells.Item(401, 13).Value = 3000
$ws.Cells.Item(401, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(401, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(401, 16).Value = 30
$ws.Cells.Item(401, 17).Value = 100
$ws.Cells.Item(401, 18).Value = "Hortaliza"
